# Commit message: "stateless entities outside the US"
#
# The workbook gains a new model variant "IMF (20%)" which is inserted
# (in the shared-string table / conceptually) right before the existing
# "IMF" variant, and the old "OECD (20%)" variant is dropped.  Each block
# of 8 data columns (one block per "M_%cit" / "M_ETR" / "M_PL" section)
# lays out the models in this fixed order:
#   GFA-Sales, GFA-Sales+Emp, IMF-Sales, IMF-Sales+Emp,
#   OECD(20%)-Sales, OECD(20%)-Sales+Emp, OECD-Sales, OECD-Sales+Emp
#
# After the edit the order becomes:
#   GFA-Sales, GFA-Sales+Emp, IMF(20%)-Sales, IMF(20%)-Sales+Emp,
#   IMF-Sales, IMF-Sales+Emp, OECD-Sales, OECD-Sales+Emp
#
# i.e. columns 3-4 of each 8-column block get brand new "IMF (20%)"
# numbers, and the old columns 3-4 (which used to hold the "IMF" numbers)
# slide right into columns 5-6 (displacing the old "OECD (20%)" numbers,
# which disappear).  Columns 1-2 and 7-8 of each block are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each data block starts at this column (1-based) within columns B..Y;
# column offsets within a block (0-based): 0=Sales,1=Sales+Emp (GFA),
# 2=Sales,3=Sales+Emp (old IMF / new IMF(20%)), 4=Sales,5=Sales+Emp
# (old OECD(20%) / new IMF), 6=Sales,7=Sales+Emp (OECD).
$blockStarts = @(2, 10, 18)   # B=2, J=10, R=18 (1-based column numbers)

# New "IMF (20%)" numeric values for each row (rows 4-8), one pair
# (Sales, Sales+Emp) per row, only defined for the first block (M_%cit);
# the M_ETR and M_PL blocks keep their previous "IMF" figures for the
# new IMF(20%) slot (no freshly computed number was supplied for them).
$newImf20 = @{
    4 = @(0.01093635869920273, 0.009072691627789734)
    5 = @(0.02973125562628548, 0.0434634358097193)
    6 = @(0.01387800792749046, 0.03124031414955012)
    7 = @(0.05508305901430718, 0.04152435081943773)
    8 = @(0.002682359336613213, 0.005695999317677495)
}

for ($r = 4; $r -le 8; $r++) {
    foreach ($blockStart in $blockStarts) {
        $colSales    = $blockStart + 2   # old "IMF - Sales" column
        $colSalesEmp = $blockStart + 3   # old "IMF - Sales + Emp" column
        $colOecdSales    = $blockStart + 4   # old "OECD (20%) - Sales" column
        $colOecdSalesEmp = $blockStart + 5   # old "OECD (20%) - Sales + Emp" column

        # Remember the pre-edit values before anything is overwritten.
        $oldSales    = $ws.Cells.Item($r, $colSales).Value2
        $oldSalesEmp = $ws.Cells.Item($r, $colSalesEmp).Value2

        # Old "IMF" numbers slide right into the old "OECD (20%)" slots,
        # which is now the "IMF" slot.
        $ws.Cells.Item($r, $colOecdSales).Value2    = $oldSales
        $ws.Cells.Item($r, $colOecdSalesEmp).Value2 = $oldSalesEmp

        # The vacated "IMF" slots become "IMF (20%)".
        if ($blockStart -eq 2) {
            $vals = $newImf20[$r]
            $ws.Cells.Item($r, $colSales).Value2    = $vals[0]
            $ws.Cells.Item($r, $colSalesEmp).Value2 = $vals[1]
        } else {
            # M_ETR / M_PL blocks: no new figure supplied, keep old value.
            $ws.Cells.Item($r, $colSales).Value2    = $oldSales
            $ws.Cells.Item($r, $colSalesEmp).Value2 = $oldSalesEmp
        }
    }
}

# Update the header row labels (row 2) for every block so the column
# headings reflect the new model ordering.
foreach ($blockStart in $blockStarts) {
    $colSales        = $blockStart + 2
    $colSalesEmp     = $blockStart + 3
    $colOecdSales    = $blockStart + 4
    $colOecdSalesEmp = $blockStart + 5

    $ws.Cells.Item(2, $colOecdSales).Value2    = "IMF - Sales"
    $ws.Cells.Item(2, $colOecdSalesEmp).Value2 = "IMF - Sales + Emp"

    $ws.Cells.Item(2, $colSales).Value2    = "IMF (20%) - Sales"
    $ws.Cells.Item(2, $colSalesEmp).Value2 = "IMF (20%) - Sales + Emp"
}
